# Update the three-digit x one-digit multiplication answers in the table.
# Cells are addressed by (row, column) rather than by Find/Replace because
# several answer strings (e.g. "329×2=658") repeat verbatim in more than
# one cell but must be replaced with *different* new values depending on
# position, so a global text search-and-replace would be ambiguous.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of row -> array of new values (1 entry per populated column, in
# left-to-right order) for every data row of the table.
$updates = @{
    1  = @("533×9=4797", "910×7=6370", "472×7=3304", "881×8=7048", "836×7=5852")
    5  = @("328×2=656",  "650×2=1300", "949×9=8541", "215×2=430",  "769×7=5383")
    10 = @("534×3=1602", "554×3=1662", "590×6=3540", "149×2=298",  "506×5=2530")
    15 = @("881×3=2643", "778×6=4668", "669×2=1338", "969×8=7752", "532×9=4788")
    20 = @("817×3=2451", "483×5=2415", "137×8=1096", "302×5=1510", "666×8=5328")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
